$wb = $excel.ActiveWorkbook

# ==== Step 1: update the "总计" (summary) sheet ====
# Insert a new leading data row for 2022-Q3 and shift the rest down by one,
# growing the table from 8 to 9 data points (A1:D8 -> A1:D9).
$ws1 = $wb.Worksheets.Item("总计")

# Give new row 9 the same look (style) as row 8 before filling it with data
$ws1.Cells.Item(8,1).Copy()
$ws1.Cells.Item(9,1).PasteSpecial(-4122)

$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q3"
$ws1.Cells.Item(2,3).Value = 15
$ws1.Cells.Item(2,4).Value = 9.77
$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "2022-Q2"
$ws1.Cells.Item(3,3).Value = 56
$ws1.Cells.Item(3,4).Value = 18.32
$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,2).Value = "2022-Q1"
$ws1.Cells.Item(4,3).Value = 20
$ws1.Cells.Item(4,4).Value = 2.79
$ws1.Cells.Item(5,1).Value = 3
$ws1.Cells.Item(5,2).Value = "2021-Q4"
$ws1.Cells.Item(5,3).Value = 30
$ws1.Cells.Item(5,4).Value = 10.06
$ws1.Cells.Item(6,1).Value = 4
$ws1.Cells.Item(6,2).Value = "2021-Q3"
$ws1.Cells.Item(6,3).Value = 62
$ws1.Cells.Item(6,4).Value = 10.04
$ws1.Cells.Item(7,1).Value = 5
$ws1.Cells.Item(7,2).Value = "2021-Q2"
$ws1.Cells.Item(7,3).Value = 68
$ws1.Cells.Item(7,4).Value = 20.99
$ws1.Cells.Item(8,1).Value = 6
$ws1.Cells.Item(8,2).Value = "2021-Q1"
$ws1.Cells.Item(8,3).Value = 61
$ws1.Cells.Item(8,4).Value = 36.83
$ws1.Cells.Item(9,1).Value = 7
$ws1.Cells.Item(9,2).Value = "2020-Q4"
$ws1.Cells.Item(9,3).Value = 69
$ws1.Cells.Item(9,4).Value = 37.72

# ==== Step 2: add the new "2022-Q3" worksheet with its fund-holding detail ====
# Duplicate the existing "2022-Q2" sheet (same column layout/formatting) and
# place the copy immediately before it, so order becomes ...,2022-Q3,2022-Q2,...
$src = $wb.Worksheets.Item("2022-Q2")
$src.Copy($src, $null)

# The duplicated sheet is named "2022-Q2 (2)"; fetch it fresh and rename it
$wsNew = $wb.Worksheets.Item("2022-Q2 (2)")
$wsNew.Name = "2022-Q3"

# Re-fetch by the new name (required for edits to stick after a rename)
$ws2 = $wb.Worksheets.Item("2022-Q3")

# The source sheet has 56 data rows (rows 2-57); our new sheet only needs 15
# data rows (rows 2-16), so drop the extra rows 17-57
$ws2.Rows("17:57").Delete()

$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "'000689"
$ws2.Cells.Item(2,3).Value = "前海开源新经济灵活配置混合A"
$ws2.Cells.Item(2,4).Value = "'98.71"
$ws2.Cells.Item(2,5).Value = "'92.74"
$ws2.Cells.Item(2,6).Value = "'7.57"
$ws2.Cells.Item(2,7).Value = "'7.4723"
$ws2.Cells.Item(2,8).Value = 1
$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,2).Value = "'013157"
$ws2.Cells.Item(3,3).Value = "前海开源新经济灵活配置混合C"
$ws2.Cells.Item(3,4).Value = "'14.05"
$ws2.Cells.Item(3,5).Value = "'92.74"
$ws2.Cells.Item(3,6).Value = "'7.57"
$ws2.Cells.Item(3,7).Value = "'1.0636"
$ws2.Cells.Item(3,8).Value = 1
$ws2.Cells.Item(4,1).Value = 2
$ws2.Cells.Item(4,2).Value = "'213003"
$ws2.Cells.Item(4,3).Value = "宝盈策略增长混合"
$ws2.Cells.Item(4,4).Value = "'10.55"
$ws2.Cells.Item(4,5).Value = "'90.74"
$ws2.Cells.Item(4,6).Value = "'6.25"
$ws2.Cells.Item(4,7).Value = "'0.6594"
$ws2.Cells.Item(4,8).Value = 6
$ws2.Cells.Item(5,1).Value = 3
$ws2.Cells.Item(5,2).Value = "'213002"
$ws2.Cells.Item(5,3).Value = "宝盈泛沿海增长混合"
$ws2.Cells.Item(5,4).Value = "'5.18"
$ws2.Cells.Item(5,5).Value = "'91.39"
$ws2.Cells.Item(5,6).Value = "'6.30"
$ws2.Cells.Item(5,7).Value = "'0.3263"
$ws2.Cells.Item(5,8).Value = 6
$ws2.Cells.Item(6,1).Value = 4
$ws2.Cells.Item(6,2).Value = "'233001"
$ws2.Cells.Item(6,3).Value = "大摩基础行业混合"
$ws2.Cells.Item(6,4).Value = "'0.70"
$ws2.Cells.Item(6,5).Value = "'78.50"
$ws2.Cells.Item(6,6).Value = "'8.79"
$ws2.Cells.Item(6,7).Value = "'0.0615"
$ws2.Cells.Item(6,8).Value = 4
$ws2.Cells.Item(7,1).Value = 5
$ws2.Cells.Item(7,2).Value = "'562880"
$ws2.Cells.Item(7,3).Value = "嘉实中证电池主题ETF"
$ws2.Cells.Item(7,4).Value = "'1.92"
$ws2.Cells.Item(7,5).Value = "'98.47"
$ws2.Cells.Item(7,6).Value = "'2.45"
$ws2.Cells.Item(7,7).Value = "'0.0470"
$ws2.Cells.Item(7,8).Value = 10
$ws2.Cells.Item(8,1).Value = 6
$ws2.Cells.Item(8,2).Value = "'020034"
$ws2.Cells.Item(8,3).Value = "国泰民安增利债券C"
$ws2.Cells.Item(8,4).Value = "'1.09"
$ws2.Cells.Item(8,5).Value = "'38.31"
$ws2.Cells.Item(8,6).Value = "'3.79"
$ws2.Cells.Item(8,7).Value = "'0.0413"
$ws2.Cells.Item(8,8).Value = 4
$ws2.Cells.Item(9,1).Value = 7
$ws2.Cells.Item(9,2).Value = "'000796"
$ws2.Cells.Item(9,3).Value = "宝盈睿丰创新灵活配置混合 - C"
$ws2.Cells.Item(9,4).Value = "'0.64"
$ws2.Cells.Item(9,5).Value = "'89.70"
$ws2.Cells.Item(9,6).Value = "'6.08"
$ws2.Cells.Item(9,7).Value = "'0.0389"
$ws2.Cells.Item(9,8).Value = 6
$ws2.Cells.Item(10,1).Value = 8
$ws2.Cells.Item(10,2).Value = "'000794"
$ws2.Cells.Item(10,3).Value = "宝盈睿丰创新灵活配置混合 - A/B"
$ws2.Cells.Item(10,4).Value = "'0.41"
$ws2.Cells.Item(10,5).Value = "'89.70"
$ws2.Cells.Item(10,6).Value = "'6.08"
$ws2.Cells.Item(10,7).Value = "'0.0249"
$ws2.Cells.Item(10,8).Value = 6
$ws2.Cells.Item(11,1).Value = 9
$ws2.Cells.Item(11,2).Value = "'020033"
$ws2.Cells.Item(11,3).Value = "国泰民安增利债券A"
$ws2.Cells.Item(11,4).Value = "'0.35"
$ws2.Cells.Item(11,5).Value = "'38.31"
$ws2.Cells.Item(11,6).Value = "'3.79"
$ws2.Cells.Item(11,7).Value = "'0.0133"
$ws2.Cells.Item(11,8).Value = 4
$ws2.Cells.Item(12,1).Value = 10
$ws2.Cells.Item(12,2).Value = "'010756"
$ws2.Cells.Item(12,3).Value = "兴华永兴混合A"
$ws2.Cells.Item(12,4).Value = "'0.23"
$ws2.Cells.Item(12,5).Value = "'94.73"
$ws2.Cells.Item(12,6).Value = "'3.88"
$ws2.Cells.Item(12,7).Value = "'0.0089"
$ws2.Cells.Item(12,8).Value = 8
$ws2.Cells.Item(13,1).Value = 11
$ws2.Cells.Item(13,2).Value = "'015998"
$ws2.Cells.Item(13,3).Value = "大成中证电池主题指数C"
$ws2.Cells.Item(13,4).Value = "'0.34"
$ws2.Cells.Item(13,5).Value = "'91.77"
$ws2.Cells.Item(13,6).Value = "'2.32"
$ws2.Cells.Item(13,7).Value = "'0.0079"
$ws2.Cells.Item(13,8).Value = 10
$ws2.Cells.Item(14,1).Value = 12
$ws2.Cells.Item(14,2).Value = "'159918"
$ws2.Cells.Item(14,3).Value = "嘉实中创400ETF"
$ws2.Cells.Item(14,4).Value = "'0.57"
$ws2.Cells.Item(14,5).Value = "'98.47"
$ws2.Cells.Item(14,6).Value = "'0.86"
$ws2.Cells.Item(14,7).Value = "'0.0049"
$ws2.Cells.Item(14,8).Value = 2
$ws2.Cells.Item(15,1).Value = 13
$ws2.Cells.Item(15,2).Value = "'015997"
$ws2.Cells.Item(15,3).Value = "大成中证电池主题指数A"
$ws2.Cells.Item(15,4).Value = "'0.13"
$ws2.Cells.Item(15,5).Value = "'91.77"
$ws2.Cells.Item(15,6).Value = "'2.32"
$ws2.Cells.Item(15,7).Value = "'0.0030"
$ws2.Cells.Item(15,8).Value = 10
$ws2.Cells.Item(16,1).Value = 14
$ws2.Cells.Item(16,2).Value = "'010757"
$ws2.Cells.Item(16,3).Value = "兴华永兴混合C"
$ws2.Cells.Item(16,4).Value = "'0.01"
$ws2.Cells.Item(16,5).Value = "'94.73"
$ws2.Cells.Item(16,6).Value = "'3.88"
$ws2.Cells.Item(16,7).Value = "'0.0004"
$ws2.Cells.Item(16,8).Value = 8
